$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Primary (default) header/footer -> header2.xml / footer2.xml
$hdrPrimary = $sec.Headers(1)
$ftrPrimary = $sec.Footers(1)

# First-page header/footer -> header1.xml / footer1.xml
$hdrFirst = $sec.Headers(2)
$ftrFirst = $sec.Footers(2)

# BTec_Logo-Orange pictures: image2.jpg -> image1.jpg
if ($hdrPrimary.Range.InlineShapes.Count -gt 0) {
    $hdrPrimary.Range.InlineShapes(1).Name = "image1.jpg"
}
if ($hdrFirst.Range.InlineShapes.Count -gt 0) {
    $hdrFirst.Range.InlineShapes(1).Name = "image1.jpg"
}

# PearsonLogo pictures: image1.png -> image2.png
if ($ftrPrimary.Range.InlineShapes.Count -gt 0) {
    $ftrPrimary.Range.InlineShapes(1).Name = "image2.png"
}
if ($ftrFirst.Range.InlineShapes.Count -gt 0) {
    $ftrFirst.Range.InlineShapes(1).Name = "image2.png"
}
